$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label the existing "orientation scheme" rows to reflect the newly
# run spiral / Gaussian-quadrature schemes (values only; indices in column A
# and the HKL-count values in C:M are unchanged). ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# --- Append the three rows that fell off the end of the original list
# (they now live at rows 17-19, continuing the same A/B/C..M pattern). ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C17:M19").Value = 1

# Match the formatting already used for the index column (A10:A16 -> bold,
# bordered, centered style) on the three new index cells.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "done"
